$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.001.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").Value = "'1.846.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("D4").Value = "'1.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.62%  '

$ws.Range("D5").Value = "'1.014"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("D6").Value = "'309.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").Value = "'0.4760"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.94%  '

$ws.Range("D8").Value = "'0.3678"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.49%  '

$ws.Range("D9").Value = "'0.07232"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.20%  '

$ws.Range("D10").Value = "'0.9310"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.66%  '

$ws.Range("D11").Value = "'19.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.40%  '

$ws.Range("D12").Value = "'0.07772"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.81%  '

$ws.Range("D13").Value = "'1.813.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.43%  '

$ws.Range("D14").Value = "'5.390"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.01%  '

$ws.Range("D15").Value = "'6.478"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.86%  '

$ws.Range("D16").Value = "'88.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.69%  '

$ws.Range("D17").Value = "'1.016"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.52%  '

$ws.Range("D18").Value = "'0.000008661"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.84%  '

$ws.Range("D19").Value = "'1.014"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.55%  '

$ws.Range("D20").Value = "'27.032.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("D21").Value = "'14.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.33%  '

$ws.Range("D22").Value = "'5.059"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.67%  '

$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").Value = "'1.926"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.26%  '

$ws.Range("D25").Value = "'152.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Value = "'18.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").Value = "'1.988"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.86%  '

$ws.Range("D28").Value = "'114.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.38%  '

$ws.Range("D29").Value = "'4.954"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.78%  '

$ws.Range("D30").Value = "'0.08864"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").Value = "'3.319"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.08%  '

$ws.Range("D32").Value = "'1.180"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.28%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = "'0.7380"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.41%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = "'4.507"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.84%  '

$ws.Range("D35").Value = "'2.660"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.08%  '

$ws.Range("D36").Value = "'1.118"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.36%  '

$ws.Range("D37").Value = "'0.01970"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.52%  '

$ws.Range("D38").Value = "'0.05256"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.72%  '

$ws.Range("D39").Value = "'2.979"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("D40").Value = "'0.5262"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.60%  '

$ws.Range("D41").Value = "'7.030"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.72%  '

$ws.Range("D42").Value = "'0.1513"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.14%  '

$ws.Range("D43").Value = "'8.276"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.53%  '

$ws.Range("D44").Value = "'10.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.12%  '

$ws.Range("D45").Value = "'0.4737"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.83%  '

$ws.Range("D47").Value = "'101.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.05%  '

$ws.Range("D48").Value = "'1.607"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.15%  '

$ws.Range("D49").Value = "'65.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.71%  '

$ws.Range("D50").Value = "'0.06059"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.31%  '

$ws.Range("D51").Value = "'0.8924"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.36%  '
